# Update gh-pages to output generated at 456a3b4
# Applies numeric "want to go" (F column) bumps across sheets, and inserts a
# new row for the "广州·德国美因茨名家管弦乐团 2025 新年音乐会" event into the
# 演出 (Performance) sheet and the 全部类型 (All types) aggregate sheet.

$wb = $excel.ActiveWorkbook

function Set-NumValue($ws, $addr, $val) {
    $ws.Range($addr).Value = $val
}

function Set-TextValue($ws, $addr, $val, $cleanStyleSource) {
    # Force text storage (avoid Excel's automatic date / number parsing),
    # then strip the leftover "@" text-format style so the cell ends up with
    # the same (default/no) style as its neighbours.
    $ws.Range($addr).NumberFormat = "@"
    $ws.Range($addr).Value = $val
    if ($cleanStyleSource) {
        $ws.Range($cleanStyleSource).Copy()
        $ws.Range($addr).PasteSpecial(-4122)
    }
}

function Clear-CellStyle($ws, $addr, $sourceAddr) {
    # Reset a cell's style/number-format to match a plain, unstyled source cell
    # (paste formats only, i.e. xlPasteFormats = -4122)
    $ws.Range($sourceAddr).Copy()
    $ws.Range($addr).PasteSpecial(-4122)
}

# ---------------------------------------------------------------------------
# Sheet "展览" (Exhibitions) - F column updates
# ---------------------------------------------------------------------------
$wsExpo = $wb.Worksheets.Item("展览")
Set-NumValue $wsExpo "F2"  290
Set-NumValue $wsExpo "F4"  247
Set-NumValue $wsExpo "F5"  7
Set-NumValue $wsExpo "F6"  2058
Set-NumValue $wsExpo "F7"  220
Set-NumValue $wsExpo "F8"  653
Set-NumValue $wsExpo "F10" 172
Set-NumValue $wsExpo "F11" 147
Set-NumValue $wsExpo "F12" 643
Set-NumValue $wsExpo "F13" 36
Set-NumValue $wsExpo "F14" 81
Set-NumValue $wsExpo "F15" 1202
Set-NumValue $wsExpo "F16" 100
Set-NumValue $wsExpo "F19" 245

# ---------------------------------------------------------------------------
# Sheet "演出" (Performances) - F column updates
# ---------------------------------------------------------------------------
$wsShow = $wb.Worksheets.Item("演出")
Set-NumValue $wsShow "F9"  116
Set-NumValue $wsShow "F12" 208

# Insert a new row 20 (shifts the existing row 20 "维也纳皇家交响乐团" down to 21)
$wsShow.Rows.Item(20).Insert()
Clear-CellStyle $wsShow "A20" "A21"

Set-TextValue $wsShow "B20" "2024-12-24" "C21"
$wsShow.Range("C20").Value = "广州·德国美因茨名家管弦乐团 2025 新年音乐会"
$wsShow.Range("D20").Value = "人民北路875号（广州市少年宫内） 广州蓓蕾剧院"
$wsShow.Range("E20").Value = "2024.12.24 19:30-12.24 21:00"
Set-NumValue $wsShow "F20" 0
Set-NumValue $wsShow "G20" 126
$wsShow.Range("H20").Value = "https://show.bilibili.com/platform/detail.html?id=93359"
$wsShow.Range("I20").Value = "//i0.hdslb.com/bfs/openplatform/202410/HaoFdo471728632672864.jpeg"

$wsShow.Range("A20").Value = 19
$wsShow.Range("A21").Value = 20

# ---------------------------------------------------------------------------
# Sheet "本地生活" (Local life) - F column updates
# ---------------------------------------------------------------------------
$wsLocal = $wb.Worksheets.Item("本地生活")
Set-NumValue $wsLocal "F2" 6319
Set-NumValue $wsLocal "F3" 784
Set-NumValue $wsLocal "F4" 1971
Set-NumValue $wsLocal "F5" 199

# ---------------------------------------------------------------------------
# Sheet "全部类型" (All types aggregate) - F column updates
# ---------------------------------------------------------------------------
$wsAll = $wb.Worksheets.Item("全部类型")
Set-NumValue $wsAll "F2"  6319
Set-NumValue $wsAll "F3"  784
Set-NumValue $wsAll "F4"  1971
Set-NumValue $wsAll "F5"  290
Set-NumValue $wsAll "F6"  199
Set-NumValue $wsAll "F12" 247
Set-NumValue $wsAll "F13" 7
Set-NumValue $wsAll "F17" 2058
Set-NumValue $wsAll "F18" 116
Set-NumValue $wsAll "F19" 220
Set-NumValue $wsAll "F22" 653
Set-NumValue $wsAll "F24" 172
Set-NumValue $wsAll "F25" 208
Set-NumValue $wsAll "F26" 147
Set-NumValue $wsAll "F27" 643
Set-NumValue $wsAll "F28" 36
Set-NumValue $wsAll "F29" 81
Set-NumValue $wsAll "F31" 1203
Set-NumValue $wsAll "F32" 100

# Row 41 becomes the German orchestra concert (was the Vienna concert)
Set-TextValue $wsAll "B41" "2024-12-24" "C41"
$wsAll.Range("C41").Value = "广州·德国美因茨名家管弦乐团 2025 新年音乐会"
$wsAll.Range("D41").Value = "人民北路875号（广州市少年宫内） 广州蓓蕾剧院"
$wsAll.Range("E41").Value = "2024.12.24 19:30-12.24 21:00"
Set-NumValue $wsAll "F41" 0
Set-NumValue $wsAll "G41" 126
$wsAll.Range("H41").Value = "https://show.bilibili.com/platform/detail.html?id=93359"
$wsAll.Range("I41").Value = "//i0.hdslb.com/bfs/openplatform/202410/HaoFdo471728632672864.jpeg"

# Insert a new row 42 for the Vienna concert (shifts the "U.M.A" row down to 43)
$wsAll.Rows.Item(42).Insert()
Clear-CellStyle $wsAll "A42" "A43"

Set-TextValue $wsAll "B42" "2024-12-29" "C43"
$wsAll.Range("C42").Value = "广州·维也纳皇家交响乐团2025新年音乐会"
$wsAll.Range("D42").Value = "人民北路696号 广州友谊剧院"
$wsAll.Range("E42").Value = "2024.12.29 20:00-12.30 21:45"
Set-NumValue $wsAll "F42" 45
Set-NumValue $wsAll "G42" 280
$wsAll.Range("H42").Value = "https://show.bilibili.com/platform/detail.html?id=89837"
$wsAll.Range("I42").Value = "//i2.hdslb.com/bfs/openplatform/202407/OzlirVhz1721882951190.jpeg"

$wsAll.Range("A41").Value = 40
$wsAll.Range("A42").Value = 41
$wsAll.Range("A43").Value = 42

# Final F43 value for the shifted "U.M.A" row (was F42 = 244, now 245)
Set-NumValue $wsAll "F43" 245

Write-Host "edit complete"
